$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 0
    6 = 0
    7 = 1
    8 = 0
    9 = 1
    10 = 0
    11 = 3
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 2
    18 = 0
    19 = 4
    20 = 1
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 2
    26 = 1
    27 = 2
    28 = 2
    29 = 3
    30 = 1
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 6
    36 = 4
    37 = 0
    38 = 3
    39 = 1
    40 = 2
    41 = 4
    42 = 1
    43 = 3
    44 = 1
    45 = 0
    46 = 6
    47 = 3
    48 = 2
    49 = 0
    50 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
